$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "27.948.72"
$ws.Cells.Item(2, 5).Value = "  -2.11%  "

$ws.Cells.Item(3, 4).Value = "1.868.98"
$ws.Cells.Item(3, 5).Value = "  -2.51%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.000"
$ws.Cells.Item(4, 5).Value = "  +0.00%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "312.04"
$ws.Cells.Item(5, 5).Value = "  -1.12%  "

$ws.Cells.Item(6, 5).Value = "  +0.00%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4985"
$ws.Cells.Item(7, 5).Value = "  -3.42%  "

$ws.Cells.Item(8, 4).NumberFormat = "@"
$ws.Cells.Item(8, 4).Value = "0.3814"
$ws.Cells.Item(8, 5).Value = "  -4.44%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.08919"
$ws.Cells.Item(9, 5).Value = "  -9.52%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "1.118"
$ws.Cells.Item(10, 5).Value = "  -2.95%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "41.45"
$ws.Cells.Item(11, 5).Value = "  -1.98%  "

$ws.Cells.Item(12, 4).NumberFormat = "@"
$ws.Cells.Item(12, 4).Value = "6.308"
$ws.Cells.Item(12, 5).Value = "  -3.29%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "20.66"
$ws.Cells.Item(13, 5).Value = "  -2.51%  "

$ws.Cells.Item(14, 4).Value = "1.860.51"
$ws.Cells.Item(14, 5).Value = "  -3.10%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "7.232"
$ws.Cells.Item(15, 5).Value = "  -3.23%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "1.001"
$ws.Cells.Item(16, 5).Value = "  +0.01%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "0.00001099"
$ws.Cells.Item(17, 5).Value = "  -3.48%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "90.83"
$ws.Cells.Item(18, 5).Value = "  -4.05%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "0.06629"
$ws.Cells.Item(19, 5).Value = "  -0.50%  "

$ws.Cells.Item(20, 4).NumberFormat = "@"
$ws.Cells.Item(20, 4).Value = "17.90"
$ws.Cells.Item(20, 5).Value = "  -1.89%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "6.100"
$ws.Cells.Item(22, 5).Value = "  -3.41%  "

$ws.Cells.Item(23, 4).Value = "27.963.48"
$ws.Cells.Item(23, 5).Value = "  -2.26%  "

$ws.Cells.Item(24, 4).NumberFormat = "@"
$ws.Cells.Item(24, 4).Value = "11.47"
$ws.Cells.Item(24, 5).Value = "  -0.38%  "

$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.283"
$ws.Cells.Item(25, 5).Value = "  -1.50%  "

$ws.Cells.Item(26, 2).Value = "LEO"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "3.384"
$ws.Cells.Item(26, 5).Value = "  +0.12%  "

$ws.Cells.Item(27, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(27, 4).Value = "2.074.10"
$ws.Cells.Item(27, 5).Value = "  -2.85%  "

$ws.Cells.Item(28, 2).Value = "LidoDAOToken"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "2.515"
$ws.Cells.Item(28, 5).Value = "  -6.28%  "

$ws.Cells.Item(29, 2).Value = "Monero"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "157.79"
$ws.Cells.Item(29, 5).Value = "  +0.12%  "

$ws.Cells.Item(30, 2).Value = "EthereumClassic"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "20.70"
$ws.Cells.Item(30, 5).Value = "  -2.93%  "

$ws.Cells.Item(31, 2).Value = "BitcoinCash"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "126.21"
$ws.Cells.Item(31, 5).Value = "  -2.78%  "

$ws.Cells.Item(32, 2).Value = "Stellar"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.1056"
$ws.Cells.Item(32, 5).Value = "  -1.74%  "

$ws.Cells.Item(33, 2).Value = "ImmutableX"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "1.056"
$ws.Cells.Item(33, 5).Value = "  -5.30%  "

$ws.Cells.Item(34, 2).Value = "Filecoin"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "5.574"
$ws.Cells.Item(34, 5).Value = "  -2.75%  "

$ws.Cells.Item(35, 2).Value = "HuobiToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.588"
$ws.Cells.Item(35, 5).Value = "  -1.19%  "

$ws.Cells.Item(36, 2).Value = "FraxShare"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "9.332"
$ws.Cells.Item(36, 5).Value = "  -5.19%  "

$ws.Cells.Item(37, 2).Value = "Hedera"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.06536"
$ws.Cells.Item(37, 5).Value = "  -3.43%  "

$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.02393"
$ws.Cells.Item(38, 5).Value = "  -1.95%  "

$ws.Cells.Item(39, 2).Value = "Algorand"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "0.2195"
$ws.Cells.Item(39, 5).Value = "  -1.59%  "

$ws.Cells.Item(40, 2).Value = "TrustWalletToken"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "1.291"
$ws.Cells.Item(40, 5).Value = "  +8.63%  "

$ws.Cells.Item(41, 2).Value = "ARBITRUM"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "1.198"
$ws.Cells.Item(41, 5).Value = "  -5.91%  "

$ws.Cells.Item(42, 2).Value = "Aptos"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "11.66"
$ws.Cells.Item(42, 5).Value = "  -1.41%  "

$ws.Cells.Item(43, 2).Value = "TheSandbox"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "0.6369"
$ws.Cells.Item(43, 5).Value = "  -1.89%  "

$ws.Cells.Item(44, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "4.898"
$ws.Cells.Item(44, 5).Value = "  -3.85%  "

$ws.Cells.Item(45, 2).Value = "Frax"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "1.0000"
$ws.Cells.Item(45, 5).Value = "  -0.01%  "

$ws.Cells.Item(46, 2).Value = "EnergySwap"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "13.16"
$ws.Cells.Item(46, 5).Value = "  -3.29%  "

$ws.Cells.Item(47, 2).Value = "Decentraland"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.5999"
$ws.Cells.Item(47, 5).Value = "  -1.78%  "

$ws.Cells.Item(48, 2).Value = "WEMIXTOKEN"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.280"
$ws.Cells.Item(48, 5).Value = "  -0.53%  "

$ws.Cells.Item(49, 2).Value = "PancakeSwap"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "3.669"
$ws.Cells.Item(49, 5).Value = "  -2.82%  "

$ws.Cells.Item(50, 2).Value = "EOS"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "1.227"
$ws.Cells.Item(50, 5).Value = "  +1.67%  "

$ws.Cells.Item(51, 2).Value = "NEARProtocol"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.975"
$ws.Cells.Item(51, 5).Value = "  -4.72%  "
